# Add a "publisher" column to the edit list, as the first data column
# (right after the journal's "update" column A). Every row in this sheet
# is a "sage" journal, so the new column is filled with that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The existing hyperlinks (row 3, currently columns C:F) will shift one
# column to the right once column B is inserted, but this engine's
# column-insert does not itself relocate hyperlink anchors/rels. Capture
# the original link targets (URL-escaped, same as the stored rels) now,
# before the shift, so they can be re-attached afterwards.
$urlAbstract = "http://pps.sagepub.com/content/%5b0-9%5d%7b1,%7d/%5b0-9%5d%7b1,%7d/%5b0-9%5d%7b1,%7d.abstract"  # -> will live at G3
$urlVolIssue = "http://pps.sagepub.com/content/vol%5b0-9%5d%7b1,%7d/issue%5b0-9%5d%7b1,%7d/"                   # -> will live at F3
$urlByYear4  = "http://pps.sagepub.com/content/by/year/%5b0-9%5d%7b4%7d"                                       # -> will live at E3
$urlByYear   = "http://pps.sagepub.com/content/by/year/"                                                       # -> will live at D3

# Drop the old hyperlinks so we don't leave stale ones behind pointing at
# the wrong (pre-shift) cells; we'll recreate them after the insert. The
# cell text itself (already the shared-string URL) is untouched by this.
$ws.Range("C3").Hyperlinks.Delete()

# Insert a new column before column B; B:F shift right to C:G.
$ws.Columns("B").Insert()

# New "publisher" column.
$ws.Range("B1").Value = "publisher"
$ws.Range("B2").Value = "sage"
$ws.Range("B3").Value = "sage"
$ws.Range("B4").Value = "sage"

# Recreate the hyperlinks at their shifted locations, adding them in the
# same order the original relationships were defined so the underlying
# rels line up with the (shifted) cells the same way they used to.
$ws.Hyperlinks.Add($ws.Range("G3"), $urlAbstract)
$ws.Hyperlinks.Add($ws.Range("F3"), $urlVolIssue)
$ws.Hyperlinks.Add($ws.Range("E3"), $urlByYear4)
$ws.Hyperlinks.Add($ws.Range("D3"), $urlByYear)

# Hyperlinks.Add reapplies the built-in Hyperlink style on its own, but
# make sure it matches the style already used by the rest of that row.
$ws.Range("D3:G3").Style = "Hyperlink"

# Leave the selection where the author ended up after typing the data.
$ws.Range("B5").Select()
